# Update the raw timing measurements on Sheet1 (new sparse-model run results)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B2").Value = 0.58799999999999997
$ws1.Range("C2").Value = 2.3769999999999998
$ws1.Range("D2").Value = 0.70099999999999996
$ws1.Range("E2").Value = 1.776
$ws1.Range("F2").Value = 5.4420000000000002

$ws1.Range("B3").Value = 0.54200000000000004
$ws1.Range("C3").Value = 2.2970000000000002
$ws1.Range("D3").Value = 0.62
$ws1.Range("E3").Value = 1.1870000000000001
$ws1.Range("F3").Value = 4.6459999999999999

$ws1.Range("B4").Value = 0.50600000000000001
$ws1.Range("C4").Value = 2.278
$ws1.Range("D4").Value = 0.63900000000000001
$ws1.Range("E4").Value = 1.5029999999999999
$ws1.Range("F4").Value = 4.9260000000000002
$ws1.Range("G4").Value = 1

$ws1.Range("B5").Value = 0.29499999999999998
$ws1.Range("C5").Value = 2.3109999999999999
$ws1.Range("D5").Value = 0.59699999999999998
$ws1.Range("E5").Value = 1.3240000000000001
$ws1.Range("F5").Value = 4.5259999999999998
$ws1.Range("H5").Value = 28

$ws1.Range("B6").Value = 0.72599999999999998
$ws1.Range("C6").Value = 3.1230000000000002
$ws1.Range("D6").Value = 0.6
$ws1.Range("E6").Value = 1.5109999999999999
$ws1.Range("F6").Value = 5.96

$ws1.Range("B7").Value = 0.34
$ws1.Range("C7").Value = 2.4049999999999998
$ws1.Range("D7").Value = 0.59799999999999998
$ws1.Range("E7").Value = 1.5589999999999999
$ws1.Range("F7").Value = 4.9020000000000001

$ws1.Range("B8").Value = 0.31
$ws1.Range("C8").Value = 2.6269999999999998
$ws1.Range("D8").Value = 0.56399999999999995
$ws1.Range("E8").Value = 1.4139999999999999
$ws1.Range("F8").Value = 4.9160000000000004
$ws1.Range("H8").Value = 35

$ws1.Range("B9").Value = 0.33700000000000002
$ws1.Range("C9").Value = 2.3210000000000002
$ws1.Range("D9").Value = 0.58399999999999996
$ws1.Range("E9").Value = 1.2270000000000001
$ws1.Range("F9").Value = 4.468

$ws1.Range("B10").Value = 0.46200000000000002
$ws1.Range("C10").Value = 2.2690000000000001
$ws1.Range("D10").Value = 0.58699999999999997
$ws1.Range("E10").Value = 1.49
$ws1.Range("F10").Value = 4.8079999999999998
$ws1.Range("G10").Value = 0.96299999999999997

$ws1.Range("B11").Value = 0.38800000000000001
$ws1.Range("C11").Value = 2.2130000000000001
$ws1.Range("D11").Value = 0.57999999999999996
$ws1.Range("E11").Value = 1.4530000000000001
$ws1.Range("F11").Value = 4.6340000000000003
$ws1.Range("G11").Value = 0.96499999999999997
$ws1.Range("H11").Value = 29

$ws1.Range("B12").Value = 0.374
$ws1.Range("C12").Value = 2.3860000000000001
$ws1.Range("D12").Value = 0.59799999999999998
$ws1.Range("E12").Value = 1.8939999999999999
$ws1.Range("F12").Value = 5.2510000000000003
$ws1.Range("G12").Value = 0.72699999999999998

$ws1.Range("B13").Value = 0.30499999999999999
$ws1.Range("C13").Value = 2.5329999999999999
$ws1.Range("D13").Value = 0.50600000000000001
$ws1.Range("E13").Value = 2.1019999999999999
$ws1.Range("F13").Value = 5.4450000000000003
$ws1.Range("G13").Value = 0.90700000000000003

$ws1.Range("B14").Value = 0.316
$ws1.Range("C14").Value = 2.3959999999999999
$ws1.Range("D14").Value = 0.621
$ws1.Range("E14").Value = 1.746
$ws1.Range("F14").Value = 5.0780000000000003
$ws1.Range("G14").Value = 0.73799999999999999

$ws1.Range("B15").Value = 0.39200000000000002
$ws1.Range("C15").Value = 2.4020000000000001
$ws1.Range("D15").Value = 0.61099999999999999
$ws1.Range("E15").Value = 1.5660000000000001
$ws1.Range("F15").Value = 4.9720000000000004
$ws1.Range("H15").Value = 31

$ws1.Range("B16").Value = 0.41899999999999998
$ws1.Range("C16").Value = 2.3660000000000001
$ws1.Range("D16").Value = 0.54800000000000004
$ws1.Range("E16").Value = 4.851
$ws1.Range("F16").Value = 8.1839999999999993

$ws1.Range("B17").Value = 0.32600000000000001
$ws1.Range("C17").Value = 2.6360000000000001
$ws1.Range("D17").Value = 0.56799999999999995
$ws1.Range("E17").Value = 2.5529999999999999
$ws1.Range("F17").Value = 6.0830000000000002
$ws1.Range("G17").Value = 0.69599999999999995

$ws1.Range("B18").Value = 0.30299999999999999
$ws1.Range("C18").Value = 2.3109999999999999
$ws1.Range("D18").Value = 0.61499999999999999
$ws1.Range("E18").Value = 1.306
$ws1.Range("F18").Value = 4.5339999999999998
$ws1.Range("H18").Value = 31

$ws1.Range("B19").Value = 0.34300000000000003
$ws1.Range("C19").Value = 2.407
$ws1.Range("D19").Value = 0.55900000000000005
$ws1.Range("E19").Value = 1.5089999999999999
$ws1.Range("F19").Value = 4.8239999999999998
$ws1.Range("H19").Value = 22

$ws1.Range("B20").Value = 0.66100000000000003
$ws1.Range("C20").Value = 2.4500000000000002
$ws1.Range("D20").Value = 0.55900000000000005
$ws1.Range("E20").Value = 1.1659999999999999
$ws1.Range("F20").Value = 4.835

$ws1.Range("B21").Value = 0.35899999999999999
$ws1.Range("C21").Value = 2.2650000000000001
$ws1.Range("D21").Value = 0.55400000000000005
$ws1.Range("E21").Value = 1.9650000000000001
$ws1.Range("F21").Value = 5.1420000000000003
$ws1.Range("H21").Value = 32


# Move the averages summary from column H to column G on the Munka1 sheet,
# and round the last (top_k) average to 2 decimals like the others.
$ws2 = $wb.Worksheets.Item("Munka1")
$ws2.Range("H5").ClearContents()
$ws2.Range("G5").Formula = "=ROUND(AVERAGE(Sheet1!B2:B21),2)"
$ws2.Range("H7").ClearContents()
$ws2.Range("G7").Formula = "=ROUND(AVERAGE(Sheet1!C2:C21),2)"
$ws2.Range("H9").ClearContents()
$ws2.Range("G9").Formula = "=ROUND(AVERAGE(Sheet1!D2:D21),2)"
$ws2.Range("H11").ClearContents()
$ws2.Range("G11").Formula = "=ROUND(AVERAGE(Sheet1!E2:E21),2)"
$ws2.Range("H13").ClearContents()
$ws2.Range("G13").Formula = "=ROUND(AVERAGE(Sheet1!F2:F21),2)"
$ws2.Range("H15").ClearContents()
$ws2.Range("G15").Formula = "=ROUND(AVERAGE(Sheet1!G2:G21),2)"
$ws2.Range("H17").ClearContents()
$ws2.Range("G17").Formula = "=ROUND(AVERAGE(Sheet1!H2:H21),2)"
$ws2.Range("G18").Select()
